$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.539.12'
$ws.Range("E2").Value = '  -1.48%  '
$ws.Range("D3").Value = '1.846.56'
$ws.Range("E3").Value = '  -1.13%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.72'
$ws.Range("E5").Value = '  -1.55%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  -0.14%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4235'
$ws.Range("E7").Value = '  -2.52%  '
$ws.Range("E8").Value = '  -3.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07272'
$ws.Range("E9").Value = '  -2.48%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8722'
$ws.Range("E10").Value = '  -6.47%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.61'
$ws.Range("E11").Value = '  -2.41%  '
$ws.Range("D12").Value = '1.863.59'
$ws.Range("E12").Value = '  -1.66%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.316'
$ws.Range("E13").Value = '  -2.03%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.496'
$ws.Range("E14").Value = '  -3.48%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06939'
$ws.Range("E15").Value = '  +1.27%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.005'
$ws.Range("E16").Value = '  +0.06%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '78.73'
$ws.Range("E17").Value = '  -2.86%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008857'
$ws.Range("E18").Value = '  -1.94%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.004'
$ws.Range("E19").Value = '  +0.04%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.34'
$ws.Range("E20").Value = '  -2.64%  '
$ws.Range("D21").Value = '27.547.86'
$ws.Range("E21").Value = '  -1.41%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.995'
$ws.Range("E22").Value = '  -2.49%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.38'
$ws.Range("E23").Value = '  -5.86%  '
$ws.Range("D24").Value = '2.079.63'
$ws.Range("E24").Value = '  -3.39%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.977'
$ws.Range("E25").Value = '  -3.12%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '153.22'
$ws.Range("E26").Value = '  +0.13%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.90'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '120.41'
$ws.Range("E28").Value = '  +6.54%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.232'
$ws.Range("E29").Value = '  -5.40%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.886'
$ws.Range("E30").Value = '  +11.61%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08897'
$ws.Range("E31").Value = '  -1.17%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7588'
$ws.Range("E32").Value = '  -5.79%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.552'
$ws.Range("E33").Value = '  -4.83%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.939'
$ws.Range("E34").Value = '  -0.27%  '
$ws.Range("E35").Value = '  -6.91%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.002'
$ws.Range("E36").Value = '  -0.16%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.089'
$ws.Range("E37").Value = '  -2.57%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05335'
$ws.Range("E38").Value = '  -2.93%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01932'
$ws.Range("E39").Value = '  -2.24%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.795'
$ws.Range("E40").Value = '  -6.84%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.880'
$ws.Range("E41").Value = '  -1.21%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5087'
$ws.Range("E42").Value = '  -3.00%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1646'
$ws.Range("E43").Value = '  -2.90%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.270'
$ws.Range("E44").Value = '  -5.71%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.06527'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4738'
$ws.Range("E46").Value = '  -2.80%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.38'
$ws.Range("E47").Value = '  -1.11%  '
$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '104.26'
$ws.Range("E48").Value = '  -2.23%  '
$ws.Range("B49").Value = 'PaxDollar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.002'
$ws.Range("E49").Value = '  -0.08%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.618'
$ws.Range("E50").Value = '  -2.94%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '64.55'
$ws.Range("E51").Value = '  -2.69%  '
